$wb = $excel.ActiveWorkbook

# The "Repayment schedule" sheet gets a new (blank) column inserted before
# column N, shifting the old N/O/P ("Late" / heading / "Outstanding")
# columns one place to the right (O/P/Q) - this supports the new
# "Variable Instalments" loan feature mentioned in the commit message.
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Columns("N").Insert()

# The newly inserted column inherits the width of the column to its left
# (column M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and leave the selection
# on cell K16, matching where the editor last clicked.
$ws.Activate()
$ws.Range("K16").Select() | Out-Null
